$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# 1. MHW sheet: update site rows 3 & 4 (new site naming LR1 / LR2)
# -------------------------------------------------------------------
$mhw = $wb.Worksheets.Item("MHW")

# Set in this order so new shared strings are registered as: LR1, LR2, mhwLR1, mhwLR2
$mhw.Range("B3").Value = "LR1"
$mhw.Range("B4").Value = "LR2"
$mhw.Range("A3").Value = "mhwLR1"
$mhw.Range("A4").Value = "mhwLR2"
$mhw.Range("D3").Value = "LR1"
$mhw.Range("D4").Value = "LR2"
$mhw.Range("E3").Value = 378787
$mhw.Range("F3").Value = 6427116
$mhw.Range("E4").Value = 379059
$mhw.Range("F4").Value = 6427057

$mhw.Activate()
$mhw.Range("A1:J4").Select()

# -------------------------------------------------------------------
# 2. New HYD2O sheet, copied from 360E so it inherits header styling,
#    then overwrite the data rows with the new Hyd2O logger data.
# -------------------------------------------------------------------
$e360 = $wb.Worksheets.Item("360E")

# Update the saved selection on the 360E sheet
$e360.Activate()
$e360.Range("C33").Select()

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$e360.Copy([System.Reflection.Missing]::Value, $lastSheet)

$hyd2o = $wb.Worksheets.Item($wb.Worksheets.Count)
$hyd2o.Name = "HYD2O"

# Clear the data rows that were copied from 360E
$hyd2o.Range("A3:K9").ClearContents()

# Fill new content in an order that keeps new shared strings grouped:
# EcoLR1..4, hyd2oEcoLR1..4, Lake Richmond Site 3/4
$hyd2o.Range("B3").Value = "EcoLR1"
$hyd2o.Range("B4").Value = "EcoLR2"
$hyd2o.Range("B5").Value = "EcoLR3"
$hyd2o.Range("B6").Value = "EcoLR4"

$hyd2o.Range("A3").Value = "hyd2oEcoLR1"
$hyd2o.Range("A4").Value = "hyd2oEcoLR2"
$hyd2o.Range("A5").Value = "hyd2oEcoLR3"
$hyd2o.Range("A6").Value = "hyd2oEcoLR4"

$hyd2o.Range("C3").Value = "Lake Richmond North"
$hyd2o.Range("C4").Value = "Lake Richmond Site 2"
$hyd2o.Range("C5").Value = "Lake Richmond Site 3"
$hyd2o.Range("C6").Value = "Lake Richmond Site 4"

$hyd2o.Range("D3").Value = "EcoLR1"
$hyd2o.Range("D4").Value = "EcoLR2"
$hyd2o.Range("D5").Value = "EcoLR3"
$hyd2o.Range("D6").Value = "EcoLR4"

$hyd2o.Range("E3").Value = 379234
$hyd2o.Range("F3").Value = 6426700
$hyd2o.Range("E4").Value = 379055
$hyd2o.Range("F4").Value = 6427391
$hyd2o.Range("E5").Value = 378703
$hyd2o.Range("F5").Value = 6427545
$hyd2o.Range("E6").Value = 378806
$hyd2o.Range("F6").Value = 6426899

$hyd2o.Range("J3").Value = "SW"
$hyd2o.Range("J4").Value = "SW"
$hyd2o.Range("J5").Value = "SW"
$hyd2o.Range("J6").Value = "SW"

# Match column A width from the original authored file
$hyd2o.Columns.Item(1).ColumnWidth = 14

# Activate HYD2O last so it becomes the selected/visible tab, with the
# saved selection on cell I7 like the authored workbook.
$hyd2o.Activate()
$hyd2o.Range("I7").Select()
